# Add : implement Excel file upload handling and data processing with multer
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells for the two added columns
$ws.Range("F1").Value = "Amount"
$ws.Range("G1").Value = "UserId"

# Row 2 - RETHEESH R
$ws.Range("A2").Value = "RETHEESH R"
$ws.Range("B2").Value = " rrmulayara@gmail.com"
$ws.Range("C2").Value = "ALCHP-INS-202509987"
$ws.Range("D2").Value = " rrmulayara@gmail.com"
$ws.Range("E2").Value = "ALCHP-INS-202509987"
$ws.Range("F2").Value = 25023
$ws.Range("G2").Value = "CROWN-109987"

# Row 3 - Ranju Kanat
$ws.Range("A3").Value = "Ranju Kanat"
$ws.Range("B3").Value = "menonranju777@gmail.com"
$ws.Range("C3").Value = "ALCHP-INS-202513805"
$ws.Range("D3").Value = "menonranju777@gmail.com"
$ws.Range("E3").Value = "ALCHP-INS-202513805"
$ws.Range("F3").Value = 26234
$ws.Range("G3").Value = "CROWN-113805"

# Row 4 - Padma Priya (new)
$ws.Range("A4").Value = "Padma Priya"
$ws.Range("B4").Value = "padmapriya.ranju@gmail.com"
$ws.Range("C4").Value = "ALCHP-INS-202514173"
$ws.Range("D4").Value = "padmapriya.ranju@gmail.com"
$ws.Range("E4").Value = "ALCHP-INS-202514173"
$ws.Range("F4").Value = 20010
$ws.Range("G4").Value = "CROWN-114173"

# Row 5 - Srikanth (new)
$ws.Range("A5").Value = "Srikanth"
$ws.Range("B5").Value = "menonranju777@gmail.com"
$ws.Range("C5").Value = "ALCHP-INS-202514734"
$ws.Range("D5").Value = "menonranju777@gmail.com"
$ws.Range("E5").Value = "ALCHP-INS-202514734"
$ws.Range("F5").Value = 28025
$ws.Range("G5").Value = "CROWN-114734"

# Row 6 - Abraham (new)
$ws.Range("A6").Value = "Abraham "
$ws.Range("B6").Value = "ebrothekkekara@gmail.com "
$ws.Range("C6").Value = "ALCHP-INS-202512447"
$ws.Range("D6").Value = "ebrothekkekara@gmail.com "
$ws.Range("E6").Value = "ALCHP-INS-202512447"
$ws.Range("F6").Value = 26212
$ws.Range("G6").Value = "CROWN-112447"

# Row 7 - Ally (new)
$ws.Range("A7").Value = "Ally"
$ws.Range("B7").Value = "ebrothekkekara@gmail.com "
$ws.Range("C7").Value = "ALCHP-INS-202513061"
$ws.Range("D7").Value = "ebrothekkekara@gmail.com "
$ws.Range("E7").Value = "ALCHP-INS-202513061"
$ws.Range("F7").Value = 26230
$ws.Range("G7").Value = "CROWN-113061"
